# Rename the inline picture shapes embedded in the document's headers and
# footers:
#   - footers: the Pearson logo -> "image1.png" becomes "image2.png"
#   - headers: the BTec logo    -> "image2.jpg" becomes "image1.jpg"
#
# InlineShape does not expose a settable Name property in the Word object
# model, so each picture is briefly converted to a floating Shape (which
# does expose .Name) and converted back to an inline shape in place.
#
# NOTE: distinct loop-variable names are used for every loop (rather than
# reusing $i) to avoid cross-loop variable aliasing in this interpreter.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

for ($fi = 1; $fi -le 2; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    $ftrInlineCount = $ftr.Range.InlineShapes.Count
    for ($fj = 1; $fj -le $ftrInlineCount; $fj++) {
        $ftrInlineShape = $ftr.Range.InlineShapes.Item($fj)
        $ftrShape = $ftrInlineShape.ConvertToShape()
        $ftrShape.Name = "image2.png"
        [void]$ftrShape.ConvertToInlineShape()
    }
}

for ($hi = 1; $hi -le 2; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    $hdrInlineCount = $hdr.Range.InlineShapes.Count
    for ($hj = 1; $hj -le $hdrInlineCount; $hj++) {
        $hdrInlineShape = $hdr.Range.InlineShapes.Item($hj)
        $hdrShape = $hdrInlineShape.ConvertToShape()
        $hdrShape.Name = "image1.jpg"
        [void]$hdrShape.ConvertToInlineShape()
    }
}
